$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spread the random-variable topics across two weeks, folding expectation
# into the discussion of each kind of random variable. This shifts the
# "Continuous random variables" topic down one row and replaces the old
# dedicated "Expectation" / "Convergence" rows with the two random-variable
# topics (Discrete, then Continuous).
$ws.Range("D9").Value = "Discrete random variables"
$ws.Range("D10").Value = "Continuous random variables"
$ws.Range("D11").Value = "Continuous random variables"

# Update the active selection to reflect where editing ended up.
$ws.Range("D12").Select()
